$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Config")

# Row 12
$ws.Range("C12").Value = 51
$ws.Range("D12").Value = 66275542
$ws.Range("E12").Value = 33724457
$ws.Range("F12").Value = 809
$ws.Range("G12").Value = 1739
$ws.Range("H12").Value = 2560

# Row 13
$ws.Range("C13").Value = 49
$ws.Range("D13").Value = 61918310
$ws.Range("E13").Value = 38076628
$ws.Range("F13").Value = 920
$ws.Range("G13").Value = 1762
$ws.Range("H13").Value = 2822

# Row 14
$ws.Range("C14").Value = 55
$ws.Range("D14").Value = 73724705
$ws.Range("E14").Value = 26196677
$ws.Range("F14").Value = 984
$ws.Range("G14").Value = 3317
$ws.Range("H14").Value = 4481

# Row 19
$ws.Range("C19").Value = 3201
$ws.Range("D19").Value = 6936
$ws.Range("E19").Value = 968

# Row 20
$ws.Range("C20").Value = 3363
$ws.Range("D20").Value = 6339
$ws.Range("E20").Value = 983

# Row 21
$ws.Range("C21").Value = 4252
$ws.Range("D21").Value = 7507
$ws.Range("E21").Value = 1067

# Row 25
$ws.Range("C25").Value = 30
$ws.Range("D25").Value = 30
$ws.Range("E25").Value = 30
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 56

# Row 26
$ws.Range("C26").Value = 32
$ws.Range("D26").Value = 32
$ws.Range("E26").Value = 31
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 57

# Row 27
$ws.Range("C27").Value = 34
$ws.Range("D27").Value = 34
$ws.Range("E27").Value = 32
$ws.Range("F27").Value = 32
$ws.Range("G27").Value = 54

# Row 31
$ws.Range("C31").Value = 39500911
$ws.Range("D31").Value = 39504824
$ws.Range("E31").Value = 40296194
$ws.Range("F31").Value = 40288859
$ws.Range("G31").Value = 73477930

# Row 32
$ws.Range("C32").Value = 44795490
$ws.Range("D32").Value = 44789117
$ws.Range("E32").Value = 43980234
$ws.Range("F32").Value = 43974826
$ws.Range("G32").Value = 79970422

# Row 33
$ws.Range("C33").Value = 48681368
$ws.Range("D33").Value = 48724360
$ws.Range("E33").Value = 46723542
$ws.Range("F33").Value = 46822499
$ws.Range("G33").Value = 78796462

# Row 37
$ws.Range("C37").Value = 60465361
$ws.Range("D37").Value = 60429029
$ws.Range("E37").Value = 59703806
$ws.Range("F37").Value = 59711141
$ws.Range("G37").Value = 26430944

# Row 38
$ws.Range("C38").Value = 55204510
$ws.Range("D38").Value = 55210883
$ws.Range("E38").Value = 55935594
$ws.Range("F38").Value = 55941002
$ws.Range("G38").Value = 20000350

# Row 39
$ws.Range("C39").Value = 51247387
$ws.Range("D39").Value = 51204395
$ws.Range("E39").Value = 53231255
$ws.Range("F39").Value = 53132298
$ws.Range("G39").Value = 21175897

# Row 45
$ws.Range("C45").Value = 3149
$ws.Range("D45").Value = 11761
$ws.Range("E45").Value = 931

# Row 46
$ws.Range("C46").Value = 3213
$ws.Range("D46").Value = 11245
$ws.Range("E46").Value = 1011

# Row 47
$ws.Range("C47").Value = 4252
$ws.Range("D47").Value = 7507
$ws.Range("E47").Value = 1067

# Row 51
$ws.Range("C51").Value = 22
$ws.Range("D51").Value = 22
$ws.Range("E51").Value = 22
$ws.Range("F51").Value = 23
$ws.Range("G51").Value = 23
$ws.Range("H51").Value = 23
$ws.Range("I51").Value = 15
$ws.Range("J51").Value = 15
$ws.Range("K51").Value = 15
$ws.Range("L51").Value = 15

# Row 52
$ws.Range("C52").Value = 22
$ws.Range("D52").Value = 22
$ws.Range("E52").Value = 22
$ws.Range("F52").Value = 20
$ws.Range("G52").Value = 20
$ws.Range("H52").Value = 20
$ws.Range("I52").Value = 15
$ws.Range("J52").Value = 15
$ws.Range("K52").Value = 15
$ws.Range("L52").Value = 15

# Row 53
$ws.Range("C53").Value = 21
$ws.Range("D53").Value = 21
$ws.Range("E53").Value = 21
$ws.Range("F53").Value = 20
$ws.Range("G53").Value = 20
$ws.Range("H53").Value = 20
$ws.Range("I53").Value = 16
$ws.Range("J53").Value = 16
$ws.Range("K53").Value = 16
$ws.Range("L53").Value = 16

# Row 57
$ws.Range("C57").Value = 29056198
$ws.Range("D57").Value = 29061733
$ws.Range("E57").Value = 29026758
$ws.Range("F57").Value = 29659950
$ws.Range("G57").Value = 29645723
$ws.Range("H57").Value = 29645123
$ws.Range("I57").Value = 19152994
$ws.Range("J57").Value = 19128708
$ws.Range("K57").Value = 19125895
$ws.Range("L57").Value = 19114603

# Row 58
$ws.Range("C58").Value = 31231341
$ws.Range("D58").Value = 31237644
$ws.Range("E58").Value = 31216018
$ws.Range("F58").Value = 28721135
$ws.Range("G58").Value = 28295847
$ws.Range("H58").Value = 28309343
$ws.Range("I58").Value = 21429032
$ws.Range("J58").Value = 21429032
$ws.Range("K58").Value = 21429032
$ws.Range("L58").Value = 21422162

# Row 59
$ws.Range("C59").Value = 28891566
$ws.Range("D59").Value = 28889266
$ws.Range("E59").Value = 28901717
$ws.Range("F59").Value = 27392614
$ws.Range("G59").Value = 27392356
$ws.Range("H59").Value = 27400767
$ws.Range("I59").Value = 21230891
$ws.Range("J59").Value = 21230891
$ws.Range("K59").Value = 21230891
$ws.Range("L59").Value = 21230891

# Row 63
$ws.Range("C63").Value = 70943802
$ws.Range("D63").Value = 70938267
$ws.Range("E63").Value = 70835662
$ws.Range("F63").Value = 70329343
$ws.Range("G63").Value = 70343570
$ws.Range("H63").Value = 70344170
$ws.Range("I63").Value = 80793872
$ws.Range("J63").Value = 80818518
$ws.Range("K63").Value = 80820971
$ws.Range("L63").Value = 80832263

# Row 64
$ws.Range("C64").Value = 68727073
$ws.Range("D64").Value = 68720770
$ws.Range("E64").Value = 68742396
$ws.Range("F64").Value = 71275662
$ws.Range("G64").Value = 71700950
$ws.Range("H64").Value = 71687454
$ws.Range("I64").Value = 78527265
$ws.Range("J64").Value = 78527265
$ws.Range("K64").Value = 78527265
$ws.Range("L64").Value = 78534135

# Row 65
$ws.Range("C65").Value = 71096707
$ws.Range("D65").Value = 71099007
$ws.Range("E65").Value = 71086556
$ws.Range("F65").Value = 72588989
$ws.Range("G65").Value = 72589247
$ws.Range("H65").Value = 72580836
$ws.Range("I65").Value = 78741131
$ws.Range("J65").Value = 78741131
$ws.Range("K65").Value = 78741131
$ws.Range("L65").Value = 78741131

$ws.Range("L65").Select()
